$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create row 12 with the same formatting as row 11 (copy formats only, not values)
$ws.Range("A11:BK11").Copy() | Out-Null
$ws.Range("A12:BK12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate the new row 12 (2021 data)
$ws.Range("A12").Value = "2021年"
$ws.Range("C12").Value = 494
$ws.Range("D12").Value = 215
$ws.Range("F12").Value = 19
$ws.Range("H12").Value = 1768
$ws.Range("I12").Value = 277
$ws.Range("J12").Value = 3280
$ws.Range("K12").Value = 17
$ws.Range("L12").Value = 1676
$ws.Range("M12").Value = 613
$ws.Range("N12").Value = 75
$ws.Range("P12").Value = 33
$ws.Range("Q12").Value = 2
$ws.Range("R12").Value = 4
$ws.Range("S12").Value = 253
$ws.Range("T12").Value = 7
$ws.Range("U12").Value = 29
$ws.Range("V12").Value = 154
$ws.Range("X12").Value = 475
$ws.Range("Y12").Value = 608
$ws.Range("Z12").Value = 688
$ws.Range("AA12").Value = 599
$ws.Range("AC12").Value = 78
$ws.Range("AD12").Value = 96
$ws.Range("AE12").Value = 102
$ws.Range("AF12").Value = 1015
$ws.Range("AG12").Value = 2644
$ws.Range("AH12").Value = 97
$ws.Range("AI12").Value = 19
$ws.Range("AJ12").Value = 14
$ws.Range("AK12").Value = 80
$ws.Range("AN12").Value = 199
$ws.Range("AO12").Value = 38
$ws.Range("AP12").Value = 151
$ws.Range("AR12").Value = 24
$ws.Range("AS12").Value = 260
$ws.Range("AU12").Value = 3
$ws.Range("AV12").Value = 22
$ws.Range("AW12").Value = 826
$ws.Range("AX12").Value = 82
$ws.Range("AY12").Value = 612
$ws.Range("AZ12").Value = 181
$ws.Range("BB12").Value = 77
$ws.Range("BC12").Value = 481
$ws.Range("BD12").Value = 130
$ws.Range("BE12").Value = 2268
$ws.Range("BG12").Value = 21000
$ws.Range("BH12").Value = 44
$ws.Range("BJ12").Value = 115
$ws.Range("BK12").Value = 56
